{"js": "async function insertAfterText(anchorText, newTexts) {\n  const body = context.document.body;\n  const results = body.search(anchorText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  let anchor = results.items[0].paragraphs.getFirst();\n  for (const t of newTexts) {\n    anchor = anchor.insertParagraph(t, Word.InsertLocation.after);\n  }\n  await context.sync();\n}\n\n// 1) \"Commit\" after \"Delete Weather controller\" (i.e. before \"Install Packages\")\nawait insertAfterText(\"Delete Weather controller\", [\"Commit\"]);\n\n// 2) \"Commit\" after \"Add connections string to appsettings.json\" (i.e. before \"Create Extensions folder\")\nawait insertAfterText(\"Add connections string to appsettings.json\", [\"Commit\"]);\n\n// 3) Five new steps after \"Add CofigureSqlContext method to ServiceExtensions class\"\nawait insertAfterText(\"Add CofigureSqlContext method to ServiceExtensions class\", [\n  \"Register SQL configuration in the ConfigureServices method in the Startup class\",\n  \"Commit\",\n  \"Add CORS configuration in the ServiceExtensions class\",\n  \"Register CORS configuration in the Configure Services method in the Startup class\",\n  \"Commit\",\n]);\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Insert-ListItemAfter($anchorText, [string[]]$newTexts) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $found = $range.Find.Execute($anchorText)\n    if (-not $found) {\n        throw \"Anchor text not found: $anchorText\"\n    }\n    $para = $range.Paragraphs(1)\n    $insertAfter = $para.Range\n    foreach ($t in $newTexts) {\n        $insertAfter.InsertParagraphAfter()\n        # Move to the paragraph we just created (immediately after $insertAfter's paragraph)\n        $newPara = $insertAfter.Paragraphs(1).Next()\n        $newPara.Range.Text = $t\n        $insertAfter = $newPara.Range\n    }\n}\n\n# 1) \"Commit\" after \"Delete Weather controller\" (i.e. before \"Install Packages\")\nInsert-ListItemAfter \"Delete Weather controller\" @(\"Commit\")\n\n# 2) \"Commit\" after \"Add connections string to appsettings.json\" (i.e. before \"Create Extensions folder\")\nInsert-ListItemAfter \"Add connections string to appsettings.json\" @(\"Commit\")\n\n# 3) Five new steps after \"Add CofigureSqlContext method to ServiceExtensions class\"\nInsert-ListItemAfter \"Add CofigureSqlContext method to ServiceExtensions class\" @(\n    \"Register SQL configuration in the ConfigureServices method in the Startup class\",\n    \"Commit\",\n    \"Add CORS configuration in the ServiceExtensions class\",\n    \"Register CORS configuration in the Configure Services method in the Startup class\",\n    \"Commit\"\n)\n"}
